$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task completion values (Tasks Update)
$ws.Range("D35").Value = 100
$ws.Range("D40").Value = 100
$ws.Range("D45").Value = 100
$ws.Range("D46").Value = 100

# Update the active cell selection to D33
[void]$ws.Activate()
[void]$ws.Range("D33").Select()
